# Commit: "Create Model from DXF: Load Features from Layers"
#
# The burn-down tracker's "current work" highlight band moves down one row:
#   - Row 16 ("Fix the now-broken Fit View function") gets its first actual
#     time entry (J16) and a live SUM formula in G16.
#   - Row 17 ("Show Flyover Data") is now finished: it gets an Actual Date
#     (H17) and logged time (J17), and drops out of the "current" highlight
#     band into the normal "done" band (style matches rows 4-16).
#   - Row 18 ("Load .dxf into model") becomes the new "current" row, picking
#     up the highlight band's formatting, an extra helper cell (D18) and an
#     explicit (no-longer-shared) start-date formula based on its own
#     duration estimate (E18).
# All of the downstream "Date Estimate" cells (F19:F34) are formulas chained
# off F18, so they ripple forward automatically once F18 recalculates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BurnDownList")

# --- Row 16: log time against "Fix the now-broken Fit View function" ---
$ws.Range("G16").Formula = "=SUM(J16:Y16)"
$ws.Range("J16").Value = 0.02013888888888889

# --- Move the "current row" formatting band down: 17 -> 18 ---
# Copy row 18's original ("current" band) formatting onto row 17 first is
# wrong order; row 17 must hand its band down to row 18 BEFORE row 17 itself
# is re-painted with row 16's ("done" band) formatting.
$ws.Range("A17:H17").Copy()
$ws.Range("A18:H18").PasteSpecial(-4122)
$ws.Range("J17:Y17").Copy()
$ws.Range("J18:Y18").PasteSpecial(-4122)

$ws.Range("A16:H16").Copy()
$ws.Range("A17:H17").PasteSpecial(-4122)
$ws.Range("J16:Y16").Copy()
$ws.Range("J17:Y17").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row 17: "Show Flyover Data" finished today, with time logged ---
$ws.Range("H17").Value = $ws.Range("F17").Value
$ws.Range("J17").Value = 0.07013888888888889

# --- Row 18: "Load .dxf into model" becomes the active row ---
$ws.Range("F18").Formula = "=F17+24*E18/8"

$ws.Range("F18").Select()
